$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp text
$ws.Range("A1").Value = "Datos actualizados a 23 de Mayo de 2020 a las 19:35"

# Row 4
$ws.Range("B4").Value = 67610
$ws.Range("D4").Value = 58654
$ws.Range("E4").Value = 8956

# Row 5
$ws.Range("B5").Value = 57114
$ws.Range("D5").Value = 50451
$ws.Range("E5").Value = 6663

# Row 6
$ws.Range("B6").Value = 18737
$ws.Range("D6").Value = 16776

# Row 7
$ws.Range("B7").Value = 16855
$ws.Range("D7").Value = 13921
$ws.Range("E7").Value = 2934

# Row 9
$ws.Range("B9").Value = 12597
$ws.Range("D9").Value = 11208
$ws.Range("E9").Value = 1389

# Row 13
$ws.Range("B13").Value = 5627
$ws.Range("D13").Value = 4769

# Row 16
$ws.Range("B16").Value = 5205
$ws.Range("D16").Value = 4696
$ws.Range("E16").Value = 509

# Row 20
$ws.Range("D20").Value = 3679
$ws.Range("E20").Value = 357

# Row 32
$ws.Range("B32").Value = 2377
$ws.Range("D32").Value = 2070

# Row 33
$ws.Range("B33").Value = 2321
$ws.Range("D33").Value = 2164
$ws.Range("E33").Value = 157
